# plant/Data diagram.xlsx - add Address and Room/Message mini-tables,
# shift the Group/Role/Group_Role tables two columns to the right (E/G/I -> G/I/K)
# and insert a new Address table in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be introduced in this exact order so that the
#     resulting sharedStrings.xml gets the same unique-string index sequence
#     as the target file (indices 49..60): Room, roomId, adminId, customerId,
#     Message, roomId(uuid), senderId, message, Address, wardId, provinceId,
#     districtId.
$ws.Range("C48").Value = "Room"
$ws.Range("E50").Value = "roomId"
$ws.Range("C51").Value = "adminId"
$ws.Range("C52").Value = "customerId"
$ws.Range("E48").Value = "Message"
$ws.Range("C50").Value = "roomId(uuid)"
$ws.Range("E51").Value = "senderId"
$ws.Range("E52").Value = "message"
$ws.Range("E4").Value = "Address"
$ws.Range("E9").Value = "wardId"
$ws.Range("E7").Value = "provinceId"
$ws.Range("E8").Value = "districtId"

# --- Shift the Group / Role / Group_Role tables two columns to the right.
# Group header+rows: was E4/E6/E7 -> now G4/G6/G7
$ws.Range("G4").Value = "Group"
$ws.Range("G6").Value = "name"
$ws.Range("G7").Value = "description"

# Role header+rows: was G4/G6/G7 -> now I4/I6/I7
$ws.Range("I4").Value = "Role"
$ws.Range("I6").Value = "url"
$ws.Range("I7").Value = "description"

# Group_Role header+rows: was I4/I6/I7 -> now K4/K6/K7
$ws.Range("K4").Value = "Group_Role"
$ws.Range("K6").Value = "role_Id"
$ws.Range("K7").Value = "group_Id"

# --- "id" rows for every table in row 5 (unchanged value, now also under K)
$ws.Range("C5").Value = "id"
$ws.Range("E5").Value = "id"
$ws.Range("G5").Value = "id"
$ws.Range("I5").Value = "id"
$ws.Range("K5").Value = "id"

# --- id rows for the new Room / Message tables
$ws.Range("C49").Value = "id"
$ws.Range("E49").Value = "id"

# --- header styling (orange fill, matches the other table headers)
$headerCells = @("E4","G4","I4","K4","C48","E48")
foreach ($addr in $headerCells) {
    $ws.Range($addr).Interior.Color = 49407
}

# --- style index "2" cells: blank spacer / explicitly-formatted id cells
$blankStyled = @("D4","G28","I28","G29","I29","C49","E49")
foreach ($addr in $blankStyled) {
    $ws.Range($addr).Interior.ColorIndex = -4142
}

# --- view state: selection / scroll position
$ws.Range("I54").Select() | Out-Null

Write-Host "done"
